$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before the "总计" (totals) sheet, with
#    the same layout/formatting as the other quarterly fund-holding sheets.
# ---------------------------------------------------------------------------
$srcSheet   = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Match the page margins used by the other sheets (0.75/0.75/1/1 in, 0.5/0.5 in
# header/footer == 54/54/72/72/36/36 points).
$newSheet.PageSetup.LeftMargin   = 54
$newSheet.PageSetup.RightMargin  = 54
$newSheet.PageSetup.TopMargin    = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Re-use the existing header (B1:H1) and index-column (A2:A5) formatting from
# an existing, identically-shaped sheet so no new cell styles are created.
$srcSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$srcSheet.Range("A2:A5").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)  # xlPasteFormats

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# code, name, scale, totalPosition, positionRatio, marketValue, rank
$fundRows = @(
    @("006679", "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇A", "14.75", "83.19", "2.80", "0.4130", 10),
    @("162719", "广发道琼斯美国石油开发与生产指数（QDII-LOF）A",                           "14.75", "83.19", "2.80", "0.4130", 10),
    @("006680", "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇C", "4.73",  "83.19", "2.80", "0.1324", 10),
    @("004243", "广发道琼斯美国石油开发与生产指数（QDII-LOF）C",                           "4.73",  "83.19", "2.80", "0.1324", 10)
)

# D:G hold numeric-looking figures that must stay text, like the source data.
$textRange = $newSheet.Range("B2:G5")
$textRange.NumberFormat = "@"

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Range("A$r").Value = $r - 2
    $newSheet.Range("B$r").Value = $row[0]
    $newSheet.Range("C$r").Value = $row[1]
    $newSheet.Range("D$r").Value = $row[2]
    $newSheet.Range("E$r").Value = $row[3]
    $newSheet.Range("F$r").Value = $row[4]
    $newSheet.Range("G$r").Value = $row[5]
    $newSheet.Range("H$r").Value = $row[6]
    $r++
}

$textRange.Style = "Normal"

# ---------------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the existing
#    rows down by one.
# ---------------------------------------------------------------------------
# Re-fetch the sheet by name: the worksheet collection shifted when the new
# sheet was inserted, so the earlier $totalSheet reference is stale.
$totalSheet = $wb.Worksheets.Item("总计")

# Give the new last row (row 7) the same formatting as the existing index
# column before rewriting every row's values.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)  # xlPasteFormats

$summaryRows = @(
    @("2022-Q1", 4,  1.09),
    @("2021-Q4", 4,  0.35),
    @("2021-Q3", 10, 2.17),
    @("2021-Q2", 7,  2.54),
    @("2021-Q1", 4,  0.67),
    @("2020-Q4", 7,  3.33)
)

$r = 2
foreach ($row in $summaryRows) {
    $totalSheet.Range("A$r").Value = $r - 2
    $totalSheet.Range("B$r").Value = $row[0]
    $totalSheet.Range("C$r").Value = $row[1]
    $totalSheet.Range("D$r").Value = $row[2]
    $r++
}

# Restore the original active sheet/selection.
$wb.Worksheets.Item("2020-Q4").Activate()
$wb.Worksheets.Item("2020-Q4").Range("A1").Select()
